$wb = $excel.ActiveWorkbook

$updates = @{
    "F4"  = 306
    "F6"  = 4591
    "F9"  = 1344
    "F10" = 891
    "F12" = 976
    "F14" = 543
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
